# Updates vm_pu.xlsx (res_bus voltage-magnitude results) for Case_0_79
# after re-running the case with the 380 kV slack/ext_grid setpoint lowered
# from 1.05 p.u. to 1.02 p.u. ("case with 380 kV done").
#
# All bus voltage-magnitude results (columns B-F, I-N) across the 24 result
# rows (rows 2-25) are refreshed with the newly computed per-unit values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    "B2" = 1.02
    "C2" = 1.033472160775461
    "D2" = 1.034309879571508
    "E2" = 1.041408002600597
    "F2" = 1.049445165721083
    "I2" = 1.02359499962809
    "J2" = 1.038596700309459
    "K2" = 1.037109722035751
    "L2" = 1.044187602891931
    "M2" = 1.052202204439152
    "N2" = 1.016435734895482
    "B3" = 1.02
    "C3" = 1.03544866148694
    "D3" = 1.036092549598916
    "E3" = 1.043218299993071
    "F3" = 1.051456291892649
    "I3" = 1.023504579208684
    "J3" = 1.04021072873558
    "K3" = 1.038698961400054
    "L3" = 1.045805899371433
    "M3" = 1.054022485445754
    "N3" = 1.017006909062117
    "B4" = 1.02
    "C4" = 1.036723925923738
    "D4" = 1.03724276409273
    "E4" = 1.044386385570991
    "F4" = 1.052754473939848
    "I4" = 1.023443335729026
    "J4" = 1.041251318080537
    "K4" = 1.039723590024914
    "L4" = 1.046849325335707
    "M4" = 1.055196793833903
    "N4" = 1.017374263863754
    "B5" = 1.02
    "C5" = 1.037259193151547
    "D5" = 1.037725545878715
    "E5" = 1.044876679306268
    "F5" = 1.053299497594486
    "I5" = 1.023416932628352
    "J5" = 1.041687892178773
    "K5" = 1.040153472794065
    "L5" = 1.047287109717779
    "M5" = 1.055689647449766
    "N5" = 1.017528171583642
    "B6" = 1.02
    "C6" = 1.037349017469049
    "D6" = 1.037806562578779
    "E6" = 1.044958957242604
    "F6" = 1.053390967180269
    "I6" = 1.023412460900047
    "J6" = 1.041761143266605
    "K6" = 1.040225601453058
    "L6" = 1.047360565054869
    "M6" = 1.05577235192405
    "N6" = 1.017553982574979
    "B7" = 1.02
    "C7" = 1.036731081514729
    "D7" = 1.037249218039028
    "E7" = 1.044392939896061
    "F7" = 1.052761759418986
    "I7" = 1.02344298551034
    "J7" = 1.041257155075224
    "K7" = 1.039729337537439
    "L7" = 1.046855178432661
    "M7" = 1.05520338258419
    "N7" = 1.017376322453185
    "B8" = 1.02
    "C8" = 1.034140901044235
    "D8" = 1.034913034895386
    "E8" = 1.042020496078486
    "F8" = 1.050125501164186
    "I8" = 1.023565007130014
    "J8" = 1.039142965960822
    "K8" = 1.037647593373324
    "L8" = 1.04473529534406
    "M8" = 1.052818121420896
    "N8" = 1.016629232106545
    "B9" = 1.02
    "C9" = 1.02954756710044
    "D9" = 1.03077027400614
    "E9" = 1.037813782196746
    "F9" = 1.045454937167117
    "I9" = 1.023759139049229
    "J9" = 1.035387577612421
    "K9" = 1.033950006735893
    "L9" = 1.040970454513999
    "M9" = 1.04858697709916
    "N9" = 1.015295372722479
    "B10" = 1.02
    "C10" = 1.026464305950678
    "D10" = 1.027989623692742
    "E10" = 1.034990439066695
    "F10" = 1.042322909519052
    "I10" = 1.023874603872166
    "J10" = 1.032862672350955
    "K10" = 1.03146408867505
    "L10" = 1.038439648565407
    "M10" = 1.045746068158165
    "N10" = 1.014394034301038
    "B11" = 1.02
    "C11" = 1.025123900402846
    "D11" = 1.026780828315563
    "E11" = 1.033763139320058
    "F11" = 1.040962044892833
    "I11" = 1.023921309042974
    "J11" = 1.031764036587678
    "K11" = 1.030382448758077
    "L11" = 1.037338555146756
    "M11" = 1.044510850272172
    "N11" = 1.014000783690352
    "B12" = 1.02
    "C12" = 1.02462518434037
    "D12" = 1.026331089032654
    "E12" = 1.033306523602899
    "F12" = 1.040455829196424
    "I12" = 1.023938164268129
    "J12" = 1.031355129511133
    "K12" = 1.029979872410065
    "L12" = 1.036928750016816
    "M12" = 1.044051245781936
    "N12" = 1.013854259460922
    "B13" = 1.02
    "C13" = 1.024732198639214
    "D13" = 1.026427593465892
    "E13" = 1.033404503200311
    "F13" = 1.04056444765997
    "I13" = 1.023934571058996
    "J13" = 1.031442879170226
    "K13" = 1.030066263295587
    "L13" = 1.037016691637134
    "M13" = 1.04414986864845
    "N13" = 1.013885710053951
    "B14" = 1.02
    "C14" = 1.025082693470842
    "D14" = 1.026743667941064
    "E14" = 1.033725410602492
    "F14" = 1.040920215997153
    "I14" = 1.023922712353276
    "J14" = 1.03173025316356
    "K14" = 1.030349188269273
    "L14" = 1.037304697185591
    "M14" = 1.044472875477545
    "N14" = 1.013988681248464
    "B15" = 1.02
    "C15" = 1.025298534092969
    "D15" = 1.026938313214992
    "E15" = 1.033923033182025
    "F15" = 1.041139319053472
    "I15" = 1.023915340507297
    "J15" = 1.031907203763831
    "K15" = 1.030523400046157
    "L15" = 1.037482038878835
    "M15" = 1.04467178525024
    "N15" = 1.014052064873853
    "B16" = 1.02
    "C16" = 1.026553149464475
    "D16" = 1.028069745193749
    "E16" = 1.035071788229111
    "F16" = 1.042413124769006
    "I16" = 1.023871434968932
    "J16" = 1.032935470802503
    "K16" = 1.031535761621134
    "L16" = 1.03851261219206
    "M16" = 1.045827936183995
    "N16" = 1.014420069899896
    "B17" = 1.02
    "C17" = 1.027338688704416
    "D17" = 1.028778172292842
    "E17" = 1.03579107652728
    "F17" = 1.043210878276292
    "I17" = 1.023843013862359
    "J17" = 1.033579030843156
    "K17" = 1.032169375539016
    "L17" = 1.039157645108924
    "M17" = 1.046551779665148
    "N17" = 1.014650110295421
    "B18" = 1.02
    "C18" = 1.027796367829581
    "D18" = 1.029190928668443
    "E18" = 1.036210166335453
    "F18" = 1.043675745109476
    "I18" = 1.023826118601784
    "J18" = 1.033953895354392
    "K18" = 1.032538449618426
    "L18" = 1.039533378106239
    "M18" = 1.04697349637915
    "N18" = 1.01478400341876
    "B19" = 1.02
    "C19" = 1.027952338491098
    "D19" = 1.029331591116166
    "E19" = 1.036352988039169
    "F19" = 1.043834177304363
    "I19" = 1.023820303822167
    "J19" = 1.034081628122467
    "K19" = 1.032664209843082
    "L19" = 1.039661408591046
    "M19" = 1.047117208714369
    "N19" = 1.014829609319113
    "B20" = 1.02
    "C20" = 1.027254461006635
    "D20" = 1.028702212232167
    "E20" = 1.035713951343408
    "F20" = 1.04312533354489
    "I20" = 1.023846096030143
    "J20" = 1.033510036188973
    "K20" = 1.032101446870572
    "L20" = 1.039088491457178
    "M20" = 1.046474168867376
    "N20" = 1.014625458730945
    "B21" = 1.02
    "C21" = 1.024979504615255
    "D21" = 1.026650612479897
    "E21" = 1.033630932052993
    "F21" = 1.040815471462718
    "I21" = 1.02392621804881
    "J21" = 1.031645651643397
    "K21" = 1.030265896380904
    "L21" = 1.037219909277214
    "M21" = 1.044377780019822
    "N21" = 1.013958371379995
    "B22" = 1.02
    "C22" = 1.023544334081193
    "D22" = 1.025356403140517
    "E22" = 1.032316949953445
    "F22" = 1.03935893215649
    "I22" = 1.023973741382756
    "J22" = 1.030468655047859
    "K22" = 1.029107131678038
    "L22" = 1.036040359360393
    "M22" = 1.043055116037392
    "N22" = 1.013536319569406
    "B23" = 1.02
    "C23" = 1.024305611183444
    "D23" = 1.026042902609284
    "E23" = 1.033013933171654
    "F23" = 1.040131482123312
    "I23" = 1.023948818288664
    "J23" = 1.031093064547596
    "K23" = 1.029721866172882
    "L23" = 1.036666114214742
    "M23" = 1.043756728090608
    "N23" = 1.013760308969049
    "B24" = 1.02
    "C24" = 1.027292521460822
    "D24" = 1.028736536725847
    "E24" = 1.035748802311415
    "F24" = 1.043163988914263
    "I24" = 1.023844704313072
    "J24" = 1.033541213488348
    "K24" = 1.032132142462995
    "L24" = 1.039119740571732
    "M24" = 1.046509239356229
    "N24" = 1.014636598593961
    "B25" = 1.02
    "C25" = 1.030738653274016
    "D25" = 1.031844499223362
    "E25" = 1.038904549512157
    "F25" = 1.046665513742671
    "I25" = 1.02371141742819
    "J25" = 1.036362103250247
    "K25" = 1.034909511630667
    "L25" = 1.049684282814682
    "M25" = 1.050680155565018
    "N25" = 1.015642309702364
}

foreach ($cellRef in $newValues.Keys) {
    $ws.Range($cellRef).Value = $newValues[$cellRef]
}
